$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")
if (-not $ws) { $ws = $wb.ActiveSheet }

# Update the "From" value for rule R30 (row 10, column C) from 18 to 1.
$ws.Range("C10").Value = 1
